$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contacts")

$ws.Range("A9").Value = "sid"
$ws.Range("B9").Value = "kmr"
$ws.Range("C9").Value = "sid1235@gmail.com"
